# Add the newly-collected keyword values to the bottom of the single-column
# "키워드" (keyword) list on Sheet1, continuing on from the existing A1:A83
# range down through A113, then re-point the "new entry" highlight from the
# old A72 row to the freshly appended A84 ("3XL") row and refresh the window
# view (zoom + scroll position) to match where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKeywords = @(
    '3XL',
    '(12M~18M)',
    '(1~11)',
    '(S(XS~M)~L(L~XL))',
    '(L~XL)',
    '(9~13)',
    '(1X~2X)',
    '(90~140)',
    '(3~13)',
    '(JL)',
    '(JM~JL)',
    '(주니어)',
    '(2~3)',
    ' 아동복',
    '(XXS~L)',
    '(1X~2X)',
    '(JM~JXL)',
    '(JM~JL)',
    '(2XL~3XL)',
    '(2XL~J3)',
    '(XL~2XL)',
    '(JXL~JXXL)',
    '(3XL~5X)',
    '(주니어)',
    '(21~FREE)',
    '(S(3-5세)~XL(12-14세))',
    '(1(S~M)~2(L~XL))',
    '(2XL~JL)',
    '(3XL~5X)',
    '(1호(1-2세)~6호(11-12세))'
)

$startRow = 84
for ($i = 0; $i -lt $newKeywords.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newKeywords[$i]
}

# The previous "new item" callout (yellow fill) lived on A72 - clear it now
# that a newer entry is the one being called out.
$ws.Range("A72").Interior.Pattern = -4142

# Highlight the newly-added standout entry (3XL) the same way earlier
# additions used to be flagged, just in orange this time around.
$ws.Range("A84").Interior.Color = 49407

# Scroll the window down to the newly-added rows and bump the zoom level,
# same as the user left it after appending the new keywords.
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.Zoom = 115
